$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Capture a clean "empty run + bold run" template from the
#    "Meta description" paragraph's bold run. We temporarily rewrite
#    its text to the title we ultimately want, then snapshot the
#    formatted run pair (leading empty run + bold run) so it can be
#    transplanted elsewhere with the correct text already in place.
# ------------------------------------------------------------------
$titleText = "Play Cash Compass for Free - Exciting Pirate-Themed Slot Game"
$titleLen = $titleText.Length

$metaPara = $d.Paragraphs.Item(2)
$metaStart = $metaPara.Range.Start
$metaBoldEnd = $metaStart + 16   # "Meta description" is 16 characters
$metaBoldRange = $d.Range($metaStart, $metaBoldEnd)
$metaBoldRange.Text = $titleText

$metaBoldRange2 = $d.Range($metaStart, $metaStart + $titleLen)
$boldTemplate = $metaBoldRange2.FormattedText

# ------------------------------------------------------------------
# 2) Create a clean placeholder paragraph (inherits no italic/bold)
#    right after the "May be overwhelming..." bullet -- i.e. right
#    before the final "Prompt:" paragraph -- and stamp the bold
#    template into it.
# ------------------------------------------------------------------
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*May be overwhelming for new players*") {
        $anchorPara = $d.Paragraphs.Item($i)
        break
    }
}
$anchorPara.Range.InsertParagraphAfter()
$newIdx = $anchorPara.Index + 1
$boldPara = $d.Paragraphs.Item($newIdx)
$boldPara.Range.Style = "Normal"
$boldPara.Range.FormattedText = $boldTemplate

# ------------------------------------------------------------------
# 3) Swap the text of the old "Prompt: ..." paragraph (still italic,
#    still has its leading empty run) for the meta-description text,
#    in place -- no paragraph split needed since the bold title now
#    lives in its own paragraph above.
# ------------------------------------------------------------------
$oldPrompt = 'Prompt: Create a cartoon-style feature image for the game "Cash Compass" that features a happy Maya warrior wearing glasses. The image should be eye-catching and bright, with plenty of colors to grab attention. The Maya warrior should be holding a compass and standing in front of a deserted island with a treasure chest nearby. The overall feel of the image should be adventurous and fun, with a hint of mystery and intrigue. '
$newDescription = "Read our review of Cash Compass by Hacksaw Gaming. Play this popular pirate-themed slot game for free with many bonus features and a maximum win of over 7,400x the bet."
$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newDescription, 2) | Out-Null

# ------------------------------------------------------------------
# 4) Remove the (now redundant / overwritten) "Meta description"
#    paragraph entirely, including its paragraph mark.
# ------------------------------------------------------------------
$metaPara2 = $d.Paragraphs.Item(2)
$metaPara2.Range.Delete()

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
